$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("dSF") updates - repulled data changed dSF values by -1 on these rows
$ws.Range("F2").Value = -3
$ws.Range("F7").Value = -8
$ws.Range("F10").Value = -2
$ws.Range("F15").Value = -3
$ws.Range("F19").Value = -2
